$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Kollau: 5h – Organisation/ " -> "Kollau: 5h – Organisation/" +
#    new run "Impressum einbinden", with the _GoBack bookmark repositioned
#    to sit exactly between the two runs (mirrors Word's own behaviour of
#    moving the hidden _GoBack bookmark to the most recent edit point).
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle1 = "Kollau: 5h – Organisation/ "
$idx1 = $text.IndexOf($needle1)
if ($idx1 -lt 0) { throw "needle1 not found" }

$rng1 = $d.Range($idx1, $idx1 + $needle1.Length)
$rng1.Text = "Kollau: 5h – Organisation/"
$splitPos1 = $idx1 + "Kollau: 5h – Organisation/".Length

$rngInsert1 = $d.Range($splitPos1, $splitPos1)
$rngInsert1.InsertAfter("Impressum einbinden")

# Re-locate the exact boundary after the insertion and drop the _GoBack
# bookmark there (Word keeps only one bookmark per name, so this also
# removes it from its previous location after "Easteregg").
$text = $d.Content.Text
$boundaryIdx = $text.IndexOf("Impressum einbinden")
$rngBoundary = $d.Range($boundaryIdx, $boundaryIdx)
$d.Bookmarks.Add("_GoBack", $rngBoundary) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Frühwirth: 5h – " (the one immediately followed by a line break) gets
#    "Itemgrafiken" appended before the break.
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle2 = "Frühwirth: 5h – " + [char]11
$idx2 = $text.IndexOf($needle2)
if ($idx2 -lt 0) { throw "needle2 not found" }
$insertPos2 = $idx2 + "Frühwirth: 5h – ".Length
$rngInsert2 = $d.Range($insertPos2, $insertPos2)
$rngInsert2.InsertAfter("Itemgrafiken")

# ---------------------------------------------------------------------------
# 3) Append "/PR-Folder" right after "Altmanninger: 5h – Zusammenstellen des
#    Impressums" (the 8.5. entry, not the earlier 1h one).
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle3 = "Altmanninger: 5h – Zusammenstellen des Impressums"
$idx3 = $text.IndexOf($needle3)
if ($idx3 -lt 0) { throw "needle3 not found" }
$endPos3 = $idx3 + $needle3.Length
$rngInsert3 = $d.Range($endPos3, $endPos3)
$rngInsert3.InsertAfter("/PR-Folder")
